# Generate Report for Handback
# Removes the "b1a0afd3-03b9-4d5c-850c-92a94d0b8941" row (row 3) from every
# sheet and refreshes the handback timestamps recorded for the remaining
# "95350963-d258-4712-ad8f-fc017fb12334" row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet1): drop row 3 (the b1a0afd3... file), keep row2.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Hyperlinks.Delete()
$wsOverview.Rows.Item(3).Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/64ce3823de3b2704ad11445e8c4f5d1f527f6126/e2e/95350963-d258-4712-ad8f-fc017fb12334.md",
    [Type]::Missing,
    [Type]::Missing,
    "95350963-d258-4712-ad8f-fc017fb12334.md"
)

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet2): refresh row2 handback datetimes, drop row3.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Range("E2").Value = "2016-03-17 14:39:42"
$wsZhCn.Range("H2").Value = "2016-03-17 14:40:01"
$wsZhCn.Rows.Item(3).Delete()
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/64ce3823de3b2704ad11445e8c4f5d1f527f6126/e2e/95350963-d258-4712-ad8f-fc017fb12334.md",
    [Type]::Missing,
    [Type]::Missing,
    "95350963-d258-4712-ad8f-fc017fb12334.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("B2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/64ce3823de3b2704ad11445e8c4f5d1f527f6126/e2e/95350963-d258-4712-ad8f-fc017fb12334.md",
    [Type]::Missing,
    [Type]::Missing,
    ".md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cea72255e46bf38bc7d4a175b7c7157f85a895e0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/95350963-d258-4712-ad8f-fc017fb12334.57a8c1aae103014c88ded2631b668fa271a16ea0.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "95350963-d258-4712-ad8f-fc017fb12334.57a8c1aae103014c88ded2631b668fa271a16ea0.zh-cn.xlf"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/974ffb6c735c0e86da821fc9a9314cfbfb54f867/e2e/95350963-d258-4712-ad8f-fc017fb12334.md",
    [Type]::Missing,
    [Type]::Missing,
    "95350963-d258-4712-ad8f-fc017fb12334.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3e1d00e3bc7320dfe20089eb6a78d659ec789750/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/95350963-d258-4712-ad8f-fc017fb12334.57a8c1aae103014c88ded2631b668fa271a16ea0.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "95350963-d258-4712-ad8f-fc017fb12334.57a8c1aae103014c88ded2631b668fa271a16ea0.zh-cn.xlf"
)

# ---------------------------------------------------------------------
# Sheet "de-de" (sheet3): refresh row2 handback datetimes, drop row3.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Range("E2").Value = "2016-03-17 14:39:46"
$wsDeDe.Range("H2").Value = "2016-03-17 14:40:12"
$wsDeDe.Rows.Item(3).Delete()
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/64ce3823de3b2704ad11445e8c4f5d1f527f6126/e2e/95350963-d258-4712-ad8f-fc017fb12334.md",
    [Type]::Missing,
    [Type]::Missing,
    "95350963-d258-4712-ad8f-fc017fb12334.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("B2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/64ce3823de3b2704ad11445e8c4f5d1f527f6126/e2e/95350963-d258-4712-ad8f-fc017fb12334.md",
    [Type]::Missing,
    [Type]::Missing,
    ".md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a69c7f1871b394d6d357b6f554ee29c553b34547/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/95350963-d258-4712-ad8f-fc017fb12334.57a8c1aae103014c88ded2631b668fa271a16ea0.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "95350963-d258-4712-ad8f-fc017fb12334.57a8c1aae103014c88ded2631b668fa271a16ea0.de-de.xlf"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/3328319a36208ab98d3437e29edb78fd5454a7c7/e2e/95350963-d258-4712-ad8f-fc017fb12334.md",
    [Type]::Missing,
    [Type]::Missing,
    "95350963-d258-4712-ad8f-fc017fb12334.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5a263533e482d78c8a3dbe82585d6caae49b61cc/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/95350963-d258-4712-ad8f-fc017fb12334.57a8c1aae103014c88ded2631b668fa271a16ea0.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "95350963-d258-4712-ad8f-fc017fb12334.57a8c1aae103014c88ded2631b668fa271a16ea0.de-de.xlf"
)
